$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanningUnits")

# The Identifier column now carries an explicit integer number format
# (previously it only had the "bestFit" text/general style). Apply this
# BEFORE writing the renumbered values so they land as numbers, not text.
$ws.Range("A2:A6").NumberFormat = "0"

# Fix bug with action initialisation not storing decision variable value:
# Planning unit Identifiers (column A) are renumbered to be consecutive
# (1,2,3,4,5 instead of 1,2,3,5,6), and the DownstreamId values (column B)
# are corrected to reference the (now consecutive) decision variable ids.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4

# Reflect the new active cell left behind by the edit.
$ws.Range("A7").Select() | Out-Null
